$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append a new log row (row 3) -----------------------------
$ws = $wb.Worksheets.Item("Logs")

$antwoord = "Geachte klant,`nBedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u wat meer informatie kunnen verstrekken over wat u precies wilt regelen?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"

$ws.Range("A3").Value = "Kun jij dit even regelen?"
$ws.Range("B3").Value = "mailmind.test@zohomail.eu"
$ws.Range("C3").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D3").Value = "Overig"
$ws.Range("E3").Value = $antwoord
$ws.Range("F3").Value = "2025-07-27 16:02:36"
$ws.Range("G3").Value = "Ja"
$ws.Range("H3").Value = "Nee"
$ws.Range("I3").Value = "Ja"
$ws.Range("J3").Value = "Ja"

# Extend the per-column conditional-formatting ranges down to the new row.
$ws.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D3"))
$ws.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G3"))
$ws.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H3"))
$ws.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I3"))
$ws.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J3"))

# --- "Dashboard" sheet: bump the "Overig" tally ----------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 2
